$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 01:22"

# Row 4 - Estados Unidos (United States)
$ws.Range("B4").Value = 737217
$ws.Range("C4").Value = 27482
$ws.Range("E4").Value = 630847
$ws.Range("G4").Value = 1778
$ws.Range("H4").Value = 38932

# Row 16 - Canada
$ws.Range("B16").Value = 33383
$ws.Range("C16").Value = 1456
$ws.Range("D16").Value = 11207
$ws.Range("E16").Value = 20706
$ws.Range("G16").Value = 160
$ws.Range("H16").Value = 1470

# Row 38 - Australia
$ws.Range("B38").Value = 6577
$ws.Range("C38").Value = 44
$ws.Range("D38").Value = 4167
$ws.Range("E38").Value = 2340
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 70

# Row 141 - Aruba
$ws.Range("F141").Value = 4
